$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 7 corresponds to 663c9331-3c7f-498c-b2ec-ee5e55b66c52.md
$wsOverview.Range("G7").Value = "2016-09-07 00:57:14"
$wsZhCn.Range("H7").Value = "2016-09-07 00:57:06"
$wsDeDe.Range("H7").Value = "2016-09-07 00:57:14"
